$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
# Row 43
$ws.Range("H43").Value = 52428
$ws.Range("I43").Value = 3606
$ws.Range("J43").Value = 101250
$ws.Range("K43").Value = 3606
$ws.Range("L43").Value = 101250
$ws.Range("M43").Value = -3537
$ws.Range("N43").Value = -101388
# Row 100
$ws.Range("H100").Value = 3999.6667
$ws.Range("I100").Value = 3999.6667
$ws.Range("K100").Value = 3999.6667
$ws.Range("M100").Value = -3458.6667
# Row 101
$ws.Range("H101").Value = 13327.3
$ws.Range("I101").Value = 9743.333000000001
$ws.Range("J101").Value = 18703.25
$ws.Range("K101").Value = 29229.999
$ws.Range("L101").Value = 56109.75
$ws.Range("M101").Value = -27607.999
$ws.Range("N101").Value = -59353.75
# Row 129
$ws.Range("H129").Value = 3301.5
$ws.Range("I129").Value = 3398.6667
$ws.Range("K129").Value = 10196.0001
$ws.Range("M129").Value = -5196.000100000001
# Row 137
$ws.Range("H137").Value = 2475.5264
$ws.Range("I137").Value = 1947.125
$ws.Range("J137").Value = 2859.818
$ws.Range("K137").Value = 5841.375
$ws.Range("L137").Value = 8579.454000000002
$ws.Range("M137").Value = -3291.375
$ws.Range("N137").Value = -13679.454
# Row 138
$ws.Range("H138").Value = 8543
$ws.Range("J138").Value = 8999
$ws.Range("L138").Value = 26997
$ws.Range("N138").Value = -37277

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 1288.3636
$ws.Range("I5").Value = 353.2857
$ws.Range("J5").Value = 2924.75
$ws.Range("K5").Value = 353.2857
$ws.Range("L5").Value = 2924.75
$ws.Range("M5").Value = -241.2857
$ws.Range("N5").Value = -3148.75
# Row 32
$ws.Range("H32").Value = 4163.4688
$ws.Range("I32").Value = 3655.9167
$ws.Range("J32").Value = 5686.125
$ws.Range("K32").Value = 3655.9167
$ws.Range("L32").Value = 5686.125
$ws.Range("M32").Value = -3368.9167
$ws.Range("N32").Value = -6260.125
# Row 102
$ws.Range("H102").Value = 2809.6
$ws.Range("I102").Value = 1013.7143
$ws.Range("J102").Value = 7000
$ws.Range("K102").Value = 1013.7143
$ws.Range("L102").Value = 7000
$ws.Range("M102").Value = 608.2857
$ws.Range("N102").Value = -10244
# Row 132
$ws.Range("H132").Value = 1540.7354
$ws.Range("I132").Value = 1094.6364
$ws.Range("J132").Value = 2358.5833
$ws.Range("K132").Value = 3283.9092
$ws.Range("L132").Value = 7075.749899999999
$ws.Range("M132").Value = -753.9092000000001
$ws.Range("N132").Value = -12135.7499

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 1288.3636
$ws.Range("I4").Value = 353.2857
$ws.Range("J4").Value = 2924.75
$ws.Range("K4").Value = 353.2857
$ws.Range("L4").Value = 2924.75
$ws.Range("M4").Value = -238.2857
$ws.Range("N4").Value = -3154.75
# Row 105
$ws.Range("H105").Value = 3119.3572
$ws.Range("I105").Value = 3288.5
$ws.Range("J105").Value = 2893.8333
$ws.Range("K105").Value = 3288.5
$ws.Range("L105").Value = 2893.8333
$ws.Range("M105").Value = -1541.5
$ws.Range("N105").Value = -6387.8333
# Row 134
$ws.Range("H134").Value = 2321.366
$ws.Range("I134").Value = 2143
$ws.Range("K134").Value = 6429
$ws.Range("M134").Value = -3894

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 55
$ws.Range("H55").Value = 4799
$ws.Range("I55").Value = 4799
$ws.Range("K55").Value = 4799
$ws.Range("M55").Value = -4484
# Row 58
$ws.Range("H58").Value = 6142.8887
$ws.Range("I58").Value = 6120.923
$ws.Range("K58").Value = 6120.923
$ws.Range("M58").Value = -5917.923
# Row 132
$ws.Range("H132").Value = 2667.3076
$ws.Range("I132").Value = 1967.7
$ws.Range("K132").Value = 5903.1
$ws.Range("M132").Value = -3373.1
# Row 134
$ws.Range("H134").Value = 3648.2
$ws.Range("I134").Value = 3701.7856
$ws.Range("K134").Value = 11105.3568
$ws.Range("M134").Value = -8570.356800000001
# Row 136
$ws.Range("H136").Value = 6142.8887
$ws.Range("I136").Value = 6120.923
$ws.Range("K136").Value = 18362.769
$ws.Range("M136").Value = -15812.769

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 11165579
$ws.Range("I4").Value = 1094386.2
$ws.Range("J4").Value = 200000450
$ws.Range("K4").Value = 3283158.6
$ws.Range("L4").Value = 600001350
$ws.Range("M4").Value = -3283046.6
$ws.Range("N4").Value = -600001574
# Row 5
$ws.Range("H5").Value = 9104.5
$ws.Range("J5").Value = 9531.053
$ws.Range("L5").Value = 28593.159
$ws.Range("N5").Value = -28817.159
# Row 11
$ws.Range("H11").Value = 8161.125
$ws.Range("I11").Value = 8038.533
$ws.Range("J11").Value = 10000
$ws.Range("K11").Value = 24115.599
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = -23975.599
$ws.Range("N11").Value = -30280
# Row 23
$ws.Range("H23").Value = 240
$ws.Range("J23").Value = 228.375
$ws.Range("L23").Value = 685.125
$ws.Range("N23").Value = -1155.125
# Row 32
$ws.Range("H32").Value = 50
$ws.Range("J32").Value = 50
$ws.Range("L32").Value = 150
$ws.Range("N32").Value = -716
# Row 46
$ws.Range("H46").Value = 3666.6667
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 9000
$ws.Range("M46").Value = -8909
# Row 68
$ws.Range("H68").Value = 2819.5789
$ws.Range("I68").Value = 1299.5
$ws.Range("J68").Value = 2998.4119
$ws.Range("K68").Value = 3898.5
$ws.Range("L68").Value = 8995.235700000001
$ws.Range("M68").Value = -3087.5
$ws.Range("N68").Value = -10617.2357
# Row 71
$ws.Range("H71").Value = 2819.5789
$ws.Range("I71").Value = 1299.5
$ws.Range("J71").Value = 2998.4119
$ws.Range("K71").Value = 11695.5
$ws.Range("L71").Value = 26985.7071
$ws.Range("M71").Value = -7639.5
$ws.Range("N71").Value = -35097.7071
# Row 121
$ws.Range("H121").Value = 1028.1875
$ws.Range("I121").Value = 691.5
$ws.Range("J121").Value = 1230.2
$ws.Range("K121").Value = 2074.5
$ws.Range("L121").Value = 3690.6
$ws.Range("M121").Value = -764.5
$ws.Range("N121").Value = -6310.6
# Row 131
$ws.Range("H131").Value = 21799.346
$ws.Range("J131").Value = 1681.2
$ws.Range("L131").Value = 5043.6
$ws.Range("N131").Value = -15123.6
# Row 132
$ws.Range("H132").Value = 3511.111
$ws.Range("I132").Value = 3160.2
$ws.Range("J132").Value = 3949.75
$ws.Range("K132").Value = 28441.8
$ws.Range("L132").Value = 35547.75
$ws.Range("M132").Value = -25911.8
$ws.Range("N132").Value = -40607.75
# Row 135
$ws.Range("H135").Value = 9104.5
$ws.Range("J135").Value = 9531.053
$ws.Range("L135").Value = 85779.477
$ws.Range("N135").Value = -90849.477

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 4024.9285
$ws.Range("I102").Value = 4152.1113
$ws.Range("K102").Value = 4152.1113
$ws.Range("M102").Value = -2530.1113
# Row 109
$ws.Range("H109").Value = 69995
$ws.Range("J109").Value = 69995
$ws.Range("L109").Value = 69995
$ws.Range("N109").Value = -72075
# Row 132
$ws.Range("H132").Value = 3177.0476
$ws.Range("I132").Value = 2733
$ws.Range("J132").Value = 4598
$ws.Range("K132").Value = 8199
$ws.Range("L132").Value = 13794
$ws.Range("M132").Value = -5669
$ws.Range("N132").Value = -18854

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3137.8
$ws.Range("I22").Value = 3674.75
$ws.Range("J22").Value = 990
$ws.Range("K22").Value = 3674.75
$ws.Range("L22").Value = 990
$ws.Range("M22").Value = -3379.75
$ws.Range("N22").Value = -1580
# Row 27
$ws.Range("H27").Value = 3137.8
$ws.Range("I27").Value = 3674.75
$ws.Range("J27").Value = 990
$ws.Range("K27").Value = 3674.75
$ws.Range("L27").Value = 990
$ws.Range("M27").Value = -3567.75
$ws.Range("N27").Value = -1204
# Row 121
$ws.Range("H121").Value = 95684
$ws.Range("J121").Value = 95684
$ws.Range("L121").Value = 95684
$ws.Range("N121").Value = -99178
# Row 136
$ws.Range("H136").Value = 9619.421
$ws.Range("J136").Value = 11577.5
$ws.Range("L136").Value = 34732.5
$ws.Range("N136").Value = -39832.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 121
$ws.Range("H121").Value = 109998
$ws.Range("J121").Value = 109998
$ws.Range("L121").Value = 109998
$ws.Range("N121").Value = -113492
# Row 136
$ws.Range("H136").Value = 3895.8235
$ws.Range("I136").Value = 3202.7273
$ws.Range("K136").Value = 9608.1819
$ws.Range("M136").Value = -7058.1819
